# Refresh the crypto price / 1h-volume columns (D, E) on the sheet with the
# latest scrape values, same as the GitHub Actions job that regenerates
# this workbook every run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price strings use '.' as both a thousands separator and decimal point
# (e.g. "62.869.60"), so they must stay plain text - never let Excel's
# automatic type detection reinterpret them as numbers. Forcing the Text
# number format for the duration of the write, then restoring the cell's
# original (unstyled) look via the built-in "Normal" style, achieves that
# without leaving any visible formatting change behind.
function Set-TextValue {
    param($Cell, $Text)
    $Cell.NumberFormat = "@"
    $Cell.Value = $Text
    $Cell.Style = "Normal"
}

# row -> (new Price, new Volume(1h)); Price is $null where the diff left
# column D unchanged for that row.
$updates = [ordered]@{
    2  = @("62.869.60", "  +0.07%  ")
    3  = @("3.439.37",  "  +0.00%  ")
    5  = @("576.50",    "  -1.04%  ")
    6  = @("146.37",    "  -0.06%  ")
    7  = @("3.440.88",  "  +0.06%  ")
    8  = @($null,       "  -0.06%  ")
    9  = @($null,       "  +0.12%  ")
    10 = @($null,       "  +1.29%  ")
    11 = @($null,       "  -0.89%  ")
    12 = @($null,       "  +2.53%  ")
    13 = @("4.027.61",  "  -0.14%  ")
    14 = @($null,       "  +2.53%  ")
    15 = @("28.85",     "  -1.78%  ")
    16 = @("3.435.71",  "  -0.13%  ")
    17 = @($null,       "  -0.88%  ")
    18 = @("62.919.27", "  +0.10%  ")
    19 = @("6.33",      "  +1.94%  ")
    20 = @($null,       "  +0.73%  ")
    21 = @("9.15",      "  -1.22%  ")
    22 = @("384.12",    "  -2.57%  ")
    23 = @("0.558",     "  -0.38%  ")
    24 = @("74.35",     "  -1.36%  ")
    25 = @($null,       "  -0.14%  ")
    26 = @("3.579.67",  "  -0.13%  ")
    27 = @($null,       "  -3.51%  ")
    28 = @($null,       "  -6.32%  ")
    29 = @($null,       "  -1.79%  ")
    30 = @($null,       "  +0.12%  ")
    31 = @($null,       "  -1.28%  ")
    32 = @("2.10",      "  -1.95%  ")
    33 = @("1.00",      "  -0.05%  ")
    34 = @("23.20",     "  -2.06%  ")
    35 = @("1.29",      "  -9.44%  ")
    36 = @("5.26",      "  -0.85%  ")
    37 = @($null,       "  -0.22%  ")
    38 = @("31.69",     "  +3.93%  ")
    39 = @($null,       "  -0.35%  ")
    40 = @("168.44",    "  +0.22%  ")
    41 = @("3.478.85",  "  +0.10%  ")
    42 = @("0.0766",    "  +0.23%  ")
    43 = @($null,       "  -0.41%  ")
    44 = @("42.33",     "  -1.29%  ")
    45 = @($null,       "  -0.58%  ")
    46 = @($null,       "  -0.34%  ")
    47 = @("4.33",      "  -3.27%  ")
    48 = @("2.564.00",  "  +1.58%  ")
    49 = @($null,       "  +4.10%  ")
    50 = @($null,       "  +1.22%  ")
    51 = @("22.55",     "  -4.04%  ")
}

foreach ($row in $updates.Keys) {
    $pair = $updates[$row]
    $priceText = $pair[0]
    $volumeText = $pair[1]

    if ($null -ne $priceText) {
        Set-TextValue $ws.Range("D$row") $priceText
    }

    Set-TextValue $ws.Range("E$row") $volumeText
}
